$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all cells we are about to update so Excel does not
# reinterpret numeric-looking strings (e.g. "1.00", "37.243.70") as numbers.
$cells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "E6", "E7", "D8", "E8", "E9", "E10", "E11", "D12",
    "E12", "D13", "E13", "D14", "E14", "D15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20",
    "E21", "D22", "E23", "E24", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "E30",
    "E31", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "E36", "E37", "D38", "E38", "E39", "E40",
    "D41", "E41", "D42", "E42", "E43", "D44", "E44", "E45", "D46", "E46", "E47", "D48", "E48", "D49",
    "E49", "D50", "E50", "D51", "E51"
)
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.243.70"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "2.047.62"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "230.75"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "56.99"
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "14.66"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "2.349.62"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").Value = "20.57"
$ws.Range("E14").Value = "  -3.19%  "
$ws.Range("D15").Value = "0.755"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "2.050.85"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").Value = "37.171.67"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "6.01"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "226.58"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("D26").Value = "9.65"
$ws.Range("E26").Value = "  +4.87%  "
$ws.Range("D27").Value = "169.56"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -4.67%  "
$ws.Range("D29").Value = "19.13"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("D33").Value = "0.0622"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").Value = "4.57"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").Value = "98.30"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "1.487.83"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "0.0948"
$ws.Range("E44").Value = "  -3.92%  "
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").Value = "16.39"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "3.94"
$ws.Range("E48").Value = "  -6.30%  "
$ws.Range("D49").Value = "7.24"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "2.235.72"
$ws.Range("E51").Value = "  -1.59%  "
